$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds header "K" (strikeouts), regenerated from Strike# to K.
# Update the computed K values for the affected rows.
$kValues = @{
    2  = 2
    3  = 1
    4  = 1
    5  = 0
    6  = 2
    7  = 2
    8  = 1
    9  = 7
    10 = 0
    11 = 0
    12 = 1
    13 = 1
    14 = 1
    17 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
